$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Row 1 (792÷8=, 933÷2=, 711÷7=, 858÷5=, 737÷6=)
Replace-Text "792÷8=" "662÷2="
Replace-Text "933÷2=" "203÷4="
Replace-Text "711÷7=" "660÷6="
Replace-Text "858÷5=" "231÷9="
Replace-Text "737÷6=" "158÷7="

# Row 2 (640÷2=, 727÷8=, 176÷5=, 767÷7=, 171÷3=)
# Net effect: 5 cells -> 5 cells with same formatting, content becomes
# [664÷8=, 658÷6=, 606÷8=, 727÷8=, 534÷6=]. Replace positionally to
# reproduce the identical resulting XML without altering cell count.
Replace-Text "640÷2=" "664÷8="
Replace-Text "727÷8=" "658÷6="
Replace-Text "176÷5=" "606÷8="
Replace-Text "767÷7=" "727÷8="
Replace-Text "171÷3=" "534÷6="

# Row 3 (693÷9=, 149÷2=, 630÷9=, 709÷5=, 866÷9=)
Replace-Text "693÷9=" "671÷3="
Replace-Text "149÷2=" "554÷7="
Replace-Text "630÷9=" "396÷3="
Replace-Text "709÷5=" "432÷4="
Replace-Text "866÷9=" "911÷9="

# Row 4 (660÷2=, 557÷8=, 810÷8=, 914÷6=, 191÷6=)
Replace-Text "660÷2=" "518÷6="
Replace-Text "557÷8=" "631÷6="
Replace-Text "810÷8=" "195÷4="
Replace-Text "914÷6=" "884÷4="
Replace-Text "191÷6=" "678÷4="

# Row 5 (987÷7=, 438÷4=, 829÷4=, 993÷3=, 877÷4=)
Replace-Text "987÷7=" "395÷3="
Replace-Text "438÷4=" "970÷9="
Replace-Text "829÷4=" "834÷6="
Replace-Text "993÷3=" "793÷2="
Replace-Text "877÷4=" "845÷8="
